# Fixed error where Unrealized Profit would be equal to Price Today when quantity is 0
$wb = $excel.ActiveWorkbook

# --- Summary sheet: AMD row (row 8) ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("C8").Value = 0
$summary.Range("D8").Value = 1465
$summary.Range("E8").Value = 146.64
$summary.Range("F8").Value = 0
$summary.Range("H8").Value = 0

# --- Transactions sheet: Michael B sells AMD (row 13) ---
$transactions = $wb.Worksheets.Item("Transactions")
$transactions.Range("F13").Value = 2
$transactions.Range("H13").Value = 310

# --- Michael B sheet: AMD row (row 2) ---
$michaelB = $wb.Worksheets.Item("Michael B")
$michaelB.Range("C2").Value = 0
$michaelB.Range("D2").Value = 1465
$michaelB.Range("E2").Value = 146.64
$michaelB.Range("F2").Value = 0
$michaelB.Range("H2").Value = 0
